# Coupling Config.xlsx - "changing order of CRMs in workflow"
#
# Data-level changes:
#  - Coupling Parameters!B4 (End Year): 2060 -> 2090
#  - Coupling Parameters!B18 (initialization_investment iteration): 3 -> 0
#  - Coupling Parameters!B44 (capacity_remuneration_mechanism): "strategic_reserve_ger" -> "none"
#  - Active selection on Coupling Parameters moves from C5 to C4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")
$ws.Activate() | Out-Null

$ws.Range("B4").Value = 2090
$ws.Range("B18").Value = 0
$ws.Range("B44").Value = "none"

$ws.Range("C4").Select() | Out-Null
